$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (changed) date column C for all data rows (2-13)
# from serial date 46070 (2026-02-17) to 46072 (2026-02-19)
$ws.Range("C2:C13").Value = 46072
